# edit.ps1 - applies the changes described by the diff:
#   1) Bumps every "datetimeFigureOut" date placeholder (slide master +
#      all slide layouts) from 12/7/2022 to 12/8/2022.
#   2) On slide 2's "Sumário" content placeholder, reorders the outline
#      bullets so "Definição do Problema" / "Justificativa" move up to
#      directly follow "Introdução" (ahead of "Objetivo Geral" /
#      "Objetivos Específicos"), and turns on "Shrink text on overflow"
#      (normAutofit) for that placeholder.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the date placeholders on the slide master and every layout.
# ---------------------------------------------------------------------
function Update-DateShape($shp) {
    $isDatePh = $false
    try {
        if ($shp.PlaceholderFormat.Type -eq 16) {
            $isDatePh = $true
        }
    } catch {
        $isDatePh = $false
    }
    if ($isDatePh -and $shp.HasTextFrame) {
        $shp.TextFrame.TextRange.Text = "12/8/2022"
    }
}

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape $master.Shapes.Item($i)
}

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape $layout.Shapes.Item($i)
    }
}

# ---------------------------------------------------------------------
# 2) Reorder the bullets in the "Sumário" slide's content placeholder
#    and enable normAutofit.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$newLines = @(
    "Introdução",
    "Definição do Problema",
    "Justificativa",
    "Objetivo Geral",
    "Objetivos Específicos ",
    "Turismo nas Praias",
    "Google Maps",
    "React Native",
    "Trabalho Proposto"
)
$tr.Text = ($newLines -join "`r")

# ppAutoSizeTextToFitShape -> <a:bodyPr><a:normAutofit/></a:bodyPr>
$shp.TextFrame.AutoSize = 2
